$d = $word.ActiveDocument

# Locate the "Padrão de codificação: " paragraph (the first of the three
# paragraphs being reverted/removed) and the paragraph immediately
# following the block (the next surviving paragraph, "Membros do Grupo:").
# Deleting the range between them removes:
#   - "Padrão de codificação: "
#   - <tab> + "https://google.github.io/styleguide/cppguide.html"
#   - the trailing empty paragraph
# while leaving the preceding centered empty paragraph (jc=center) intact.
$count = $d.Paragraphs.Count
$blockStartIndex = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Padrão de codificação*") {
        $blockStartIndex = $i
        break
    }
}

if ($blockStartIndex -eq -1) {
    Write-Host "Could not locate 'Padrão de codificação' paragraph; nothing to do."
} else {
    $blockEndIndex = $blockStartIndex + 2
    $start = $d.Paragraphs.Item($blockStartIndex).Range.Start
    $end = $d.Paragraphs.Item($blockEndIndex).Range.End
    $r = $d.Range($start, $end)
    $r.Delete()
    Write-Host "Removed coding-standard paragraphs."
}

# Remove the now-unused character styles (Hyperlink / Unresolved Mention)
# that were introduced alongside the hyperlink paragraph. Delete from the
# last one added back to the first, so earlier name lookups stay valid.
$d.Styles.Item("Unresolved Mention").Delete()
$d.Styles.Item("Hyperlink").Delete()
Write-Host "Removed Hyperlink / Unresolved Mention styles."
